# Updated cryptos list values (prices / 1h volume %) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be misread as a number by the
# automatic type-inference on Range.Value (e.g. "144.00" -> 144), losing
# the trailing zero(s) that the source site displays. Force text storage
# for just these cells, then restore the default "Normal" style so no
# extra formatting is left behind.
$forceTextCells = @("D6", "D12", "D17", "D29", "D33", "D36")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.960.78'
$ws.Range('E2').Value = '  -0.90%  '

$ws.Range('D3').Value = '2.504.37'
$ws.Range('E3').Value = '  +1.40%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').Value = '540.57'
$ws.Range('E5').Value = '  +0.86%  '

$ws.Range('D6').Value = '144.00'
$ws.Range('E6').Value = '  -2.57%  '

$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').Value = '0.573'
$ws.Range('E8').Value = '  +0.68%  '

$ws.Range('D9').Value = '2.528.73'
$ws.Range('E9').Value = '  +2.51%  '

$ws.Range('D10').Value = '0.101'
$ws.Range('E10').Value = '  +1.13%  '

$ws.Range('E11').Value = '  +0.42%  '

$ws.Range('D12').Value = '5.60'
$ws.Range('E12').Value = '  +5.01%  '

$ws.Range('E13').Value = '  +1.11%  '

$ws.Range('D14').Value = '2.947.27'
$ws.Range('E14').Value = '  +1.89%  '

$ws.Range('D15').Value = '23.61'
$ws.Range('E15').Value = '  -1.72%  '

$ws.Range('D16').Value = '58.890.00'
$ws.Range('E16').Value = '  -0.87%  '

$ws.Range('D17').Value = '0.0000140'
$ws.Range('E17').Value = '  +1.36%  '

$ws.Range('D18').Value = '2.521.07'
$ws.Range('E18').Value = '  +0.45%  '

$ws.Range('D19').Value = '11.22'

$ws.Range('D20').Value = '4.29'
$ws.Range('E20').Value = '  -1.49%  '

$ws.Range('D21').Value = '325.39'
$ws.Range('E21').Value = '  +0.53%  '

$ws.Range('E22').Value = '  +3.09%  '

$ws.Range('D23').Value = '5.79'
$ws.Range('E23').Value = '  +0.91%  '

$ws.Range('D24').Value = '62.03'
$ws.Range('E24').Value = '  +2.27%  '

$ws.Range('D25').Value = '0.441'
$ws.Range('E25').Value = '  -4.45%  '

$ws.Range('D26').Value = '0.163'
$ws.Range('E26').Value = '  +1.04%  '

$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.619.88'
$ws.Range('E27').Value = '  +2.21%  '

$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  +1.92%  '

$ws.Range('D29').Value = '7.80'
$ws.Range('E29').Value = '  +0.99%  '

$ws.Range('D30').Value = '0.0₃0776'
$ws.Range('E30').Value = '  +0.36%  '

$ws.Range('E31').Value = '  -0.25%  '

$ws.Range('D32').Value = '6.69'
$ws.Range('E32').Value = '  -1.35%  '

$ws.Range('D33').Value = '1.20'
$ws.Range('E33').Value = '  -4.70%  '

$ws.Range('E34').Value = '  -0.14%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '1.45'
$ws.Range('E35').Value = '  +3.78%  '

$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '156.50'
$ws.Range('E36').Value = '  +0.79%  '

$ws.Range('E37').Value = '  +1.53%  '

$ws.Range('D38').Value = '4.35'
$ws.Range('E38').Value = '  -4.27%  '

$ws.Range('D39').Value = '1.59'
$ws.Range('E39').Value = '  -8.12%  '

$ws.Range('D40').Value = '5.72'
$ws.Range('E40').Value = '  -2.63%  '

$ws.Range('D41').Value = '36.93'
$ws.Range('E41').Value = '  +0.55%  '

$ws.Range('D42').Value = '296.08'
$ws.Range('E42').Value = '  -6.28%  '

$ws.Range('E43').Value = '  -0.48%  '

$ws.Range('D44').Value = '0.822'
$ws.Range('E44').Value = '  -2.10%  '

$ws.Range('D45').Value = '0.994'
$ws.Range('E45').Value = '  -0.49%  '

$ws.Range('D46').Value = '0.601'
$ws.Range('E46').Value = '  +2.70%  '

$ws.Range('D47').Value = '10.79'
$ws.Range('E47').Value = '  +0.66%  '

$ws.Range('D48').Value = '0.0931'
$ws.Range('E48').Value = '  -0.83%  '

$ws.Range('D49').Value = '123.22'
$ws.Range('E49').Value = '  +1.37%  '

$ws.Range('D50').Value = '18.58'
$ws.Range('E50').Value = '  +0.35%  '

$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0228'
$ws.Range('E51').Value = '  -0.10%  '

# Restore the plain/default style on the forced-text cells (matches the
# original workbook, which applies no explicit number format to these).
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
